# Updated cryptos list on Fri Mar 10 08:52:51 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: D2, E2
$ws.Range("D2").Value = "'19.970.96"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -7.99%  "

# Row 3: D3, E3
$ws.Range("D3").Value = "'1.408.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -8.24%  "

# Row 4: D4
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"

# Row 5: D5, E5
$ws.Range("D5").Value = "'1.002"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.07%  "

# Row 6: D6, E6
$ws.Range("D6").Value = "'273.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.43%  "

# Row 7: D7, E7
$ws.Range("D7").Value = "'0.3700"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -6.24%  "

# Row 8: D8, E8
$ws.Range("D8").Value = "'0.3065"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.97%  "

# Row 9: D9, E9
$ws.Range("D9").Value = "'39.10"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -7.32%  "

# Row 10: B10, C10, D10, E10
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.06536"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.91%  "

# Row 11: B11, C11, D11, E11
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").Value = "'0.9887"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.52%  "

# Row 12: D12, E12
$ws.Range("D12").Value = "'1.003"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.27%  "

# Row 13: D13, E13
$ws.Range("D13").Value = "'5.292"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.82%  "

# Row 14: D14, E14
$ws.Range("D14").Value = "'6.123"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.22%  "

# Row 15: D15, E15
$ws.Range("D15").Value = "'16.83"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -8.75%  "

# Row 16: D16, E16
$ws.Range("D16").Value = "'1.410.59"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -8.06%  "

# Row 17: E17
$ws.Range("E17").Value = "  -8.50%  "

# Row 18: D18, E18
$ws.Range("D18").Value = "'0.05746"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -12.87%  "

# Row 19: D19, E19
$ws.Range("D19").Value = "'73.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -12.07%  "

# Row 20: D20, E20
$ws.Range("D20").Value = "'1.002"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.11%  "

# Row 21: D21, E21
$ws.Range("D21").Value = "'5.563"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -9.05%  "

# Row 22: D22, E22
$ws.Range("D22").Value = "'14.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.98%  "

# Row 23: E23
$ws.Range("E23").Value = "  -0.68%  "

# Row 24: D24, E24
$ws.Range("D24").Value = "'2.257"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.36%  "

# Row 25: D25, E25
$ws.Range("D25").Value = "'19.965.19"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.98%  "

# Row 26: D26, E26
$ws.Range("D26").Value = "'2.225"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.29%  "

# Row 27: D27, E27
$ws.Range("D27").Value = "'137.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.04%  "

# Row 28: D28, E28
$ws.Range("D28").Value = "'16.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.85%  "

# Row 29: D29, E29
$ws.Range("D29").Value = "'1.568.04"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.14%  "

# Row 30: D30, E30
$ws.Range("D30").Value = "'108.32"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.37%  "

# Row 31: D31, E31
$ws.Range("D31").Value = "'3.830"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -20.73%  "

# Row 32: D32, E32
$ws.Range("D32").Value = "'5.247"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -10.25%  "

# Row 33: D33, E33
$ws.Range("D33").Value = "'0.8066"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -16.00%  "

# Row 34: D34, E34
$ws.Range("D34").Value = "'0.07662"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.19%  "

# Row 35: D35, E35
$ws.Range("D35").Value = "'8.404"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.66%  "

# Row 36: D36, E36
$ws.Range("D36").Value = "'0.05788"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.51%  "

# Row 37: E37
$ws.Range("E37").Value = "  +0.07%  "

# Row 38: D38, E38
$ws.Range("D38").Value = "'4.743"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.84%  "

# Row 39: D39, E39
$ws.Range("D39").Value = "'0.1937"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.97%  "

# Row 40: D40, E40
$ws.Range("D40").Value = "'0.02028"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.64%  "

# Row 41: E41
$ws.Range("E41").Value = "  -4.50%  "

# Row 42: B42, C42, D42, E42
$ws.Range("B42").Value = "WEMIXTOKEN"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "'1.284"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -10.61%  "

# Row 43: B43, C43, D43, E43
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'1.055"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -10.38%  "

# Row 44: E44
$ws.Range("E44").Value = "  -7.93%  "

# Row 45: D45, E45
$ws.Range("D45").Value = "'3.516"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.78%  "

# Row 46: D46, E46
$ws.Range("D46").Value = "'12.04"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.53%  "

# Row 47: D47, E47
$ws.Range("D47").Value = "'0.5089"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.89%  "

# Row 48: E48
$ws.Range("E48").Value = "  -3.61%  "

# Row 49: D49, E49
$ws.Range("D49").Value = "'110.20"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.87%  "

# Row 50: D50, E50
$ws.Range("D50").Value = "'1.031"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -10.87%  "

# Row 51: D51, E51
$ws.Range("D51").Value = "'1.002"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.13%  "
